$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "time_taken" metadata
$ws.Range("F1").Value = "time_taken"

# Give the new header cell (F1) the same formatting as the other header
# cells (bold, bordered, centered) by copying the format from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# time_taken values for rows 2-27 (matching source data, microsecond precision)
$times = @(
    "2021-10-05 10:51:56.799697",
    "2021-10-05 10:51:56.799710",
    "2021-10-05 10:51:56.799714",
    "2021-10-05 10:51:56.799717",
    "2021-10-05 10:51:56.799721",
    "2021-10-05 10:51:56.799724",
    "2021-10-05 10:51:56.799727",
    "2021-10-05 10:51:56.799730",
    "2021-10-05 10:51:56.799733",
    "2021-10-05 10:51:56.799736",
    "2021-10-05 10:51:56.799739",
    "2021-10-05 10:51:56.799743",
    "2021-10-05 10:51:56.799746",
    "2021-10-05 10:51:56.799748",
    "2021-10-05 10:51:56.799751",
    "2021-10-05 10:51:56.799754",
    "2021-10-05 10:51:56.799758",
    "2021-10-05 10:51:56.799761",
    "2021-10-05 10:51:56.799764",
    "2021-10-05 10:51:56.799767",
    "2021-10-05 10:51:56.799770",
    "2021-10-05 10:51:56.799773",
    "2021-10-05 10:51:56.799776",
    "2021-10-05 10:51:56.799780",
    "2021-10-05 10:51:56.799787",
    "2021-10-05 10:51:56.799790"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
